# Harmonize similar tags to be the same (Tags section on SwateTemplateMetadata sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")
$ws.Activate()

# C12 ("assay protocol") -> "assay"
$ws.Range("C12").Value = "assay"

# C13 ("http://purl.obolibrary.org/obo/DPBO_1000177") -> "OBI:0000070"
$ws.Range("C13").Value = "OBI:0000070"

# C14 ("DPBO") -> cleared entirely (no longer needed)
$ws.Range("C14").ClearContents()

# Update active selection to reflect the new focus cell
$ws.Range("C14").Select()
